$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.0478857529451
$ws.Range("D2").Value = 1.046512829503352
$ws.Range("E2").Value = 1.061276380398277
$ws.Range("F2").Value = 1.068238907166386
$ws.Range("I2").Value = 1.035728032744373
$ws.Range("J2").Value = 1.052932439659501
$ws.Range("K2").Value = 1.049278056197399
$ws.Range("L2").Value = 1.064000880456622
$ws.Range("M2").Value = 1.070944623112024
$ws.Range("N2").Value = 1.021287279631239
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049448347775954
$ws.Range("D3").Value = 1.04768170327464
$ws.Range("E3").Value = 1.062868270282589
$ws.Range("F3").Value = 1.069977612019127
$ws.Range("I3").Value = 1.036039315098271
$ws.Range("J3").Value = 1.054140678585356
$ws.Range("K3").Value = 1.050257659721965
$ws.Range("L3").Value = 1.065405401542126
$ws.Range("M3").Value = 1.072496975759474
$ws.Range("N3").Value = 1.021705038629623
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050457510434784
$ws.Range("D4").Value = 1.048436088384153
$ws.Range("E4").Value = 1.063896767006803
$ws.Range("F4").Value = 1.071101303633836
$ws.Range("I4").Value = 1.036238550932537
$ws.Range("J4").Value = 1.054920186214632
$ws.Range("K4").Value = 1.050889005022878
$ws.Range("L4").Value = 1.066312172081463
$ws.Range("M4").Value = 1.073499625743957
$ws.Range("N4").Value = 1.02197417359593
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050881307732991
$ws.Range("D5").Value = 1.048752770408511
$ws.Range("E5").Value = 1.064328782511696
$ws.Range("F5").Value = 1.071573386974676
$ws.Range("I5").Value = 1.036321788929191
$ws.Range("J5").Value = 1.055247347310665
$ws.Range("K5").Value = 1.051153824450363
$ws.Range("L5").Value = 1.066692896724109
$ws.Range("M5").Value = 1.073920712272654
$ws.Range("N5").Value = 1.022087037198297
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.05095243878829
$ws.Range("D6").Value = 1.04880591583792
$ws.Range("E6").Value = 1.064401298716475
$ws.Range("F6").Value = 1.071652633584421
$ws.Range("I6").Value = 1.036335734483593
$ws.Range("J6").Value = 1.055302247423145
$ws.Range("K6").Value = 1.051198253896993
$ws.Range("L6").Value = 1.066756794072681
$ws.Range("M6").Value = 1.073991389761701
$ws.Range("N6").Value = 1.022105971102989
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050463175009385
$ws.Range("D7").Value = 1.048440321710264
$ws.Range("E7").Value = 1.063902541039658
$ws.Range("F7").Value = 1.071107612865567
$ws.Range("I7").Value = 1.036239665206131
$ws.Range("J7").Value = 1.05492455988504
$ws.Range("K7").Value = 1.050892545895882
$ws.Range("L7").Value = 1.066317261224738
$ws.Range("M7").Value = 1.073505253990783
$ws.Range("N7").Value = 1.021975682784956
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048414247006311
$ws.Range("D8").Value = 1.046908264405991
$ws.Range("E8").Value = 1.061814696321726
$ws.Range("F8").Value = 1.06882679997886
$ws.Range("I8").Value = 1.035833685791641
$ws.Range("J8").Value = 1.053341251264892
$ws.Range("K8").Value = 1.049609643988774
$ws.Range("L8").Value = 1.064475974364328
$ws.Range("M8").Value = 1.071469631847487
$ws.Range("N8").Value = 1.021428709698782
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044788464451212
$ws.Range("D9").Value = 1.044193326645336
$ws.Range("E9").Value = 1.05812325109661
$ws.Range("F9").Value = 1.064796773426671
$ws.Range("I9").Value = 1.035101463986061
$ws.Range("J9").Value = 1.050533272539688
$ws.Range("K9").Value = 1.047329406364253
$ws.Range("L9").Value = 1.061215309621195
$ws.Range("M9").Value = 1.067868187924848
$ws.Range("N9").Value = 1.020455700480756
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042360356176612
$ws.Range("D10").Value = 1.042372703269392
$ws.Range("E10").Value = 1.055653350187913
$ws.Range("F10").Value = 1.062102064124562
$ws.Range("I10").Value = 1.034601853154197
$ws.Range("J10").Value = 1.048648711762142
$ws.Range("K10").Value = 1.045795686381914
$ws.Range("L10").Value = 1.059030168589279
$ws.Range("M10").Value = 1.065456923984815
$ws.Range("N10").Value = 1.019800702391624
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041306228285794
$ws.Range("D11").Value = 1.041581732242625
$ws.Range("E11").Value = 1.054581601051428
$ws.Range("F11").Value = 1.060933172569801
$ws.Range("I11").Value = 1.03438276463294
$ws.Range("J11").Value = 1.047829585567103
$ws.Range("K11").Value = 1.045128265983046
$ws.Range("L11").Value = 1.058081159320817
$ws.Range("M11").Value = 1.064410240769522
$ws.Range("N11").Value = 1.019515545156723
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040914254294083
$ws.Range("D12").Value = 1.041287527414982
$ws.Range("E12").Value = 1.054183153592566
$ws.Range("F12").Value = 1.060498671076113
$ws.Range("I12").Value = 1.03430096880316
$ws.Range("J12").Value = 1.047524850653656
$ws.Range("K12").Value = 1.044879851595461
$ws.Range("L12").Value = 1.057728219298644
$ws.Range("M12").Value = 1.064021055101405
$ws.Range("N12").Value = 1.019409390773295
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040998353427245
$ws.Range("D13").Value = 1.041350653713437
$ws.Range("E13").Value = 1.054268638047897
$ws.Range("F13").Value = 1.060591887952515
$ws.Range("I13").Value = 1.034318533184827
$ws.Range("J13").Value = 1.04759023897231
$ws.Range("K13").Value = 1.044933160286627
$ws.Range("L13").Value = 1.057803946022355
$ws.Range("M13").Value = 1.064104555081315
$ws.Range("N13").Value = 1.019432171901295
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041273836312698
$ws.Range("D14").Value = 1.041557421446714
$ws.Range("E14").Value = 1.054548672511232
$ws.Range("F14").Value = 1.060897263207703
$ws.Range("I14").Value = 1.034376011879384
$ws.Range("J14").Value = 1.047804405828052
$ws.Range("K14").Value = 1.045107742312454
$ws.Range("L14").Value = 1.058051994139593
$ws.Range("M14").Value = 1.064378078808039
$ws.Range("N14").Value = 1.019506775199829
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04144351386638
$ws.Range("D15").Value = 1.041684764233813
$ws.Range("E15").Value = 1.054721163976309
$ws.Range("F15").Value = 1.061085371774817
$ws.Range("I15").Value = 1.034411371110143
$ws.Range("J15").Value = 1.047936297875232
$ws.Range("K15").Value = 1.04521524095214
$ws.Range("L15").Value = 1.05820476671062
$ws.Range("M15").Value = 1.064546552323114
$ws.Range("N15").Value = 1.019552709614972
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042430256388605
$ws.Range("D16").Value = 1.0424251413403
$ws.Range("E16").Value = 1.055724429914935
$ws.Range("F16").Value = 1.062179595006267
$ws.Range("I16").Value = 1.034616335074612
$ws.Range("J16").Value = 1.048703008379214
$ws.Range("K16").Value = 1.045839910507996
$ws.Range("L16").Value = 1.059093090695952
$ws.Range("M16").Value = 1.065526333329459
$ws.Range("N16").Value = 1.019819594653835
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043048471516821
$ws.Range("D17").Value = 1.042888850848439
$ws.Range("E17").Value = 1.056353137296796
$ws.Range("F17").Value = 1.062865410470792
$ws.Range("I17").Value = 1.034744164277509
$ws.Range("J17").Value = 1.049183109344621
$ws.Range("K17").Value = 1.046230857880278
$ws.Range("L17").Value = 1.059649547971674
$ws.Range("M17").Value = 1.066140222353186
$ws.Range("N17").Value = 1.019986590436712
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043408802303533
$ws.Range("D18").Value = 1.043159071633731
$ws.Range("E18").Value = 1.056719633860997
$ws.Range("F18").Value = 1.063265236700011
$ws.Range("I18").Value = 1.034818459389896
$ws.Range("J18").Value = 1.049462845751445
$ws.Range("K18").Value = 1.046458572058975
$ws.Range("L18").Value = 1.059973847880102
$ws.Range("M18").Value = 1.066498044338205
$ws.Range("N18").Value = 1.020083847997485
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043531621438328
$ws.Range("D19").Value = 1.043251167243857
$ws.Range("E19").Value = 1.056844563125483
$ws.Range("F19").Value = 1.063401533861251
$ws.Range("I19").Value = 1.034843747156598
$ws.Range("J19").Value = 1.0495581783121
$ws.Range("K19").Value = 1.046536162859714
$ws.Range("L19").Value = 1.060084379920368
$ws.Range("M19").Value = 1.06662001049219
$ws.Range("N19").Value = 1.02011698525944
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042982170241552
$ws.Range("D20").Value = 1.042839125446415
$ws.Range("E20").Value = 1.056285705559042
$ws.Range("F20").Value = 1.062791849569519
$ws.Range("I20").Value = 1.034730476896454
$ws.Range("J20").Value = 1.049131629983755
$ws.Range("K20").Value = 1.046188945957436
$ws.Range("L20").Value = 1.059589873605428
$ws.Range("M20").Value = 1.066074383721267
$ws.Range("N20").Value = 1.019968688726767
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041192725270659
$ws.Range("D21").Value = 1.04149654469125
$ws.Range("E21").Value = 1.054466219155976
$ws.Range("F21").Value = 1.060807346830646
$ws.Range("I21").Value = 1.034359097359889
$ws.Range("J21").Value = 1.047741352214456
$ws.Range("K21").Value = 1.045056346251067
$ws.Range("L21").Value = 1.057978962266346
$ws.Range("M21").Value = 1.064297544050517
$ws.Range("N21").Value = 1.019484812880542
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04006517120186
$ws.Range("D22").Value = 1.040650075410497
$ws.Range("E22").Value = 1.053320192761652
$ws.Range("F22").Value = 1.059557735196877
$ws.Range("I22").Value = 1.034123184547982
$ws.Range("J22").Value = 1.046864476636397
$ws.Range("K22").Value = 1.044341311847021
$ws.Range("L22").Value = 1.056963592009574
$ws.Range("M22").Value = 1.063178048913907
$ws.Range("N22").Value = 1.019179223838065
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040663146118738
$ws.Range("D23").Value = 1.041099028852247
$ws.Range("E23").Value = 1.053927920633271
$ws.Range("F23").Value = 1.060220360139152
$ws.Range("I23").Value = 1.034248475969442
$ws.Range("J23").Value = 1.04732958906
$ws.Range("K23").Value = 1.044720644804726
$ws.Range("L23").Value = 1.05750210207275
$ws.Range("M23").Value = 1.06377173886564
$ws.Range("N23").Value = 1.019341352055138
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043012129748623
$ws.Range("D24").Value = 1.042861594997475
$ws.Range("E24").Value = 1.056316175731893
$ws.Range("F24").Value = 1.062825089187306
$ws.Range("I24").Value = 1.034736662455216
$ws.Range("J24").Value = 1.049154892213031
$ws.Range("K24").Value = 1.046207885136492
$ws.Range("L24").Value = 1.05961683872406
$ws.Range("M24").Value = 1.066104134134805
$ws.Range("N24").Value = 1.019976778197071
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045727694598597
$ws.Range("D25").Value = 1.044897052816554
$ws.Range("E25").Value = 1.059079108929531
$ws.Range("F25").Value = 1.065839997054736
$ws.Range("I25").Value = 1.035292770604128
$ws.Range("J25").Value = 1.051261383186651
$ws.Range("K25").Value = 1.047921264442376
$ws.Range("L25").Value = 1.062060231938421
$ws.Range("M25").Value = 1.068801018312616
$ws.Range("N25").Value = 1.020708349876284
